$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Topic column for the newly scheduled days
$ws.Range("C61").Value = "Sports Week"
$ws.Range("C62").Value = "ED Holiday"
$ws.Range("C63").Value = "Website landing page assignment"
$ws.Range("C64").Value = "Saturday Holiday"
$ws.Range("C65").Value = "Landing Page"
$ws.Range("C66").Value = "Inernal Links HTML and Character Entities"
$ws.Range("C67").Value = "CSS Table and Text Properties"
$ws.Range("C68").Value = "Ram Nawami Holiday"

# Adjust row heights for rows 60 through 77 (slightly taller), and row 78 (tallest)
for ($r = 60; $r -le 77; $r++) {
    $ws.Rows.Item($r).RowHeight = 19.5
}
$ws.Rows.Item(78).RowHeight = 20.25
